# Auto-generated: update market-price-derived columns (H..N) across all class sheets
# per the scheduled-runner refresh. No formulas are involved anywhere in this
# workbook - every cell in H:N is a static value, so we just overwrite / clear cells.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 178.42857
$ws.Range("I9").Value = 158.16667
$ws.Range("K9").Value = 158.16667
$ws.Range("M9").Value = 10.83332999999999

$ws.Range("H12").Value = 847.7222
$ws.Range("I12").Value = 81.38461
$ws.Range("K12").Value = 81.38461
$ws.Range("M12").Value = 88.61539

$ws.Range("H15").Value = 2127.818
$ws.Range("I15").Value = 2127.818
$ws.Range("K15").Value = 6383.454000000001
$ws.Range("M15").Value = -6214.454000000001

$ws.Range("H53").Value = 1629.625
$ws.Range("I53").Value = 291.85715
$ws.Range("J53").Value = 2670.111
$ws.Range("K53").Value = 291.85715
$ws.Range("L53").Value = 2670.111
$ws.Range("M53").Value = 345.14285
$ws.Range("N53").Value = -3944.111

$ws.Range("H69").Value = 8338.75
$ws.Range("J69").Value = 8338.75
$ws.Range("L69").Value = 25016.25
$ws.Range("N69").Value = -26764.25

$ws.Range("H72").Value = 8338.75
$ws.Range("J72").Value = 8338.75
$ws.Range("L72").Value = 75048.75
$ws.Range("N72").Value = -83784.75

$ws.Range("H112").Value = 1849.0769
$ws.Range("I112").Value = 488
$ws.Range("J112").Value = 3015.7144
$ws.Range("K112").Value = 1464
$ws.Range("L112").Value = 9047.143199999999
$ws.Range("M112").Value = -356
$ws.Range("N112").Value = -11263.1432

$ws.Range("H132").Value = 1182.3784
$ws.Range("I132").Value = 1090.1765
$ws.Range("K132").Value = 3270.5295
$ws.Range("M132").Value = -740.5295000000001

$ws.Range("H138").Value = 2600.8408
$ws.Range("I138").Value = 1949.6818
$ws.Range("J138").Value = 3252
$ws.Range("K138").Value = 5849.0454
$ws.Range("L138").Value = 9756
$ws.Range("M138").Value = -709.0454
$ws.Range("N138").Value = -20036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7708.25
$ws.Range("I61").Value = 4712.5
$ws.Range("K61").Value = 4712.5
$ws.Range("M61").Value = -4500.5

$ws.Range("H74").Value = 41670708
$ws.Range("I74").Value = 66671932
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 66671932
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -66671058
$ws.Range("N74").Value = -3748

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H77").Value = 41670708
$ws.Range("I77").Value = 66671932
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 333359660
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -333355292
$ws.Range("N77").Value = -18736

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H122").Value = 4424.625
$ws.Range("I122").Value = 2224.5
$ws.Range("J122").Value = 6624.75
$ws.Range("K122").Value = 6673.5
$ws.Range("L122").Value = 19874.25
$ws.Range("M122").Value = -4223.5
$ws.Range("N122").Value = -24774.25

$ws.Range("H136").Value = 7708.25
$ws.Range("I136").Value = 4712.5
$ws.Range("K136").Value = 14137.5
$ws.Range("M136").Value = -11587.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1629.75
$ws.Range("I99").Value = 1629.75
$ws.Range("K99").Value = 1629.75
$ws.Range("M99").Value = -131.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 133.25
$ws.Range("I7").Value = 54.6
$ws.Range("J7").Value = 264.33334
$ws.Range("K7").Value = 54.6
$ws.Range("L7").Value = 264.33334
$ws.Range("M7").Value = 58.4
$ws.Range("N7").Value = -490.33334

$ws.Range("H31").Value = 50300.434
$ws.Range("I31").Value = 5717.364
$ws.Range("J31").Value = 91168.25
$ws.Range("K31").Value = 5717.364
$ws.Range("L31").Value = 91168.25
$ws.Range("M31").Value = -5422.364
$ws.Range("N31").Value = -91758.25

$ws.Range("H34").Value = 50300.434
$ws.Range("I34").Value = 5717.364
$ws.Range("J34").Value = 91168.25
$ws.Range("K34").Value = 5717.364
$ws.Range("L34").Value = 91168.25
$ws.Range("M34").Value = -5515.364
$ws.Range("N34").Value = -91572.25

$ws.Range("H122").Value = 5092.933
$ws.Range("I122").Value = 1782.4667
$ws.Range("J122").Value = 8403.4
$ws.Range("K122").Value = 5347.4001
$ws.Range("L122").Value = 25210.2
$ws.Range("M122").Value = -2897.4001
$ws.Range("N122").Value = -30110.2

$ws.Range("H134").Value = 3818.8667
$ws.Range("I134").Value = 2574.4443
$ws.Range("J134").Value = 5685.5
$ws.Range("K134").Value = 7723.3329
$ws.Range("L134").Value = 17056.5
$ws.Range("M134").Value = -5188.3329
$ws.Range("N134").Value = -22126.5

$ws.Range("H135").Value = 69443.11
$ws.Range("J135").Value = 69443.11
$ws.Range("L135").Value = 69443.11
$ws.Range("N135").Value = -79583.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3114064.2
$ws.Range("I4").Value = 3888980
$ws.Range("K4").Value = 11666940
$ws.Range("M4").Value = -11666828

$ws.Range("H9").Value = 366515.06
$ws.Range("I9").Value = 577156.7
$ws.Range("J9").Value = 2679.5454
$ws.Range("K9").Value = 1731470.1
$ws.Range("L9").Value = 8038.6362
$ws.Range("M9").Value = -1731246.1
$ws.Range("N9").Value = -8486.636200000001

$ws.Range("H98").Value = 1000
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 3000
$ws.Range("N98").Value = -5996

$ws.Range("H107").Value = 817.2963
$ws.Range("J107").Value = 1168.2307
$ws.Range("L107").Value = 3504.6921
$ws.Range("N107").Value = -7344.6921

$ws.Range("H113").Value = 1128.8
$ws.Range("I113").Value = 709.44446
$ws.Range("J113").Value = 1757.8334
$ws.Range("K113").Value = 2128.33338
$ws.Range("L113").Value = 5273.5002
$ws.Range("M113").Value = 41.66661999999997
$ws.Range("N113").Value = -9613.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1111.1818
$ws.Range("I2").Value = 28.571428
$ws.Range("J2").Value = 3005.75
$ws.Range("K2").Value = 28.571428
$ws.Range("L2").Value = 3005.75
$ws.Range("M2").Value = 84.428572
$ws.Range("N2").Value = -3231.75

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H97").Value = 1034.76
$ws.Range("I97").Value = 685.26666
$ws.Range("K97").Value = 685.26666
$ws.Range("M97").Value = -189.26666

$ws.Range("H122").Value = 3641.2354
$ws.Range("I122").Value = 2992.3572
$ws.Range("J122").Value = 6669.3335
$ws.Range("K122").Value = 8977.071599999999
$ws.Range("L122").Value = 20008.0005
$ws.Range("M122").Value = -6527.071599999999
$ws.Range("N122").Value = -24908.0005

$ws.Range("H126").Value = 3759.7727
$ws.Range("I126").Value = 2988.9443
$ws.Range("K126").Value = 8966.832900000001
$ws.Range("M126").Value = -6496.832900000001

$ws.Range("H132").Value = 64283.117
$ws.Range("I132").Value = 103175.7
$ws.Range("K132").Value = 309527.1
$ws.Range("M132").Value = -306997.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4414.231
$ws.Range("I22").Value = 2362.375
$ws.Range("J22").Value = 7697.2
$ws.Range("K22").Value = 2362.375
$ws.Range("L22").Value = 7697.2
$ws.Range("M22").Value = -2067.375
$ws.Range("N22").Value = -8287.200000000001

$ws.Range("H27").Value = 4414.231
$ws.Range("I27").Value = 2362.375
$ws.Range("J27").Value = 7697.2
$ws.Range("K27").Value = 2362.375
$ws.Range("L27").Value = 7697.2
$ws.Range("M27").Value = -2255.375
$ws.Range("N27").Value = -7911.2

$ws.Range("H46").Value = 2460.5386
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 2453.3635
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 2453.3635
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -2829.3635

$ws.Range("H55").Value = 2502470.5
$ws.Range("I55").Value = 5000994
$ws.Range("K55").Value = 5000994
$ws.Range("M55").Value = -5000821

$ws.Range("H100").Value = 9889.333000000001
$ws.Range("J100").Value = 19004
$ws.Range("L100").Value = 19004
$ws.Range("N100").Value = -20086

$ws.Range("H122").Value = 7276.5264
$ws.Range("I122").Value = 6149.6
$ws.Range("J122").Value = 11502.5
$ws.Range("K122").Value = 18448.8
$ws.Range("L122").Value = 34507.5
$ws.Range("M122").Value = -15998.8
$ws.Range("N122").Value = -39407.5

$ws.Range("H132").Value = 3075.1292
$ws.Range("I132").Value = 1318.3158
$ws.Range("J132").Value = 5856.75
$ws.Range("K132").Value = 3954.9474
$ws.Range("L132").Value = 17570.25
$ws.Range("M132").Value = -1424.9474
$ws.Range("N132").Value = -22630.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5698.346
$ws.Range("I132").Value = 6082.7144
$ws.Range("K132").Value = 18248.1432
$ws.Range("M132").Value = -15718.1432
